# Week 5 result update:
#   1. Footer "last edited" date field on the slide master + every slide
#      layout rolls from 2025-04-01 to 2025-04-02.
#   2. The two "예상 운용 가능 회비" (expected usable dues) figures on
#      slide 4 are updated: 1,889,597 -> 1,943,597 and 1,648,328 -> 1,702,328.

$p = $ppt.ActivePresentation

# --- 1. Roll the cached datetimeFigureOut field text forward one day ---
function Update-DateField($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "2025-04-01") {
                $shp.TextFrame.TextRange.Text = "2025-04-02"
            }
        }
    }
}

Update-DateField($p.SlideMaster)
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    Update-DateField($p.SlideMaster.CustomLayouts.Item($L))
}

# --- 2. Update the two due-amount figures on slide 4 ---
$slide4 = $p.Slides.Item(4)
$tr = $slide4.Shapes.Item(1).TextFrame.TextRange

$f1 = $tr.Find(": 1,889,597 ")
$f1.Text = ": 1,943,597 "

$f2 = $tr.Find(": 1,648,328 ")
$f2.Text = ": 1,702,328 "
